$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the vuln_rule identifiers in column A (rows 2-9) with "V-"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $cell.Value2 = "V-" + $val
    }
}
